# Add a new weekly record at row 194 ("Fruta / hortaliza, semanal").
# This shifts all existing rows from 194..259 down to 195..260 and we
# populate the freed-up row 194 with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 194; everything below (194-259) moves
# down to (195-260). The worksheet's used range grows from R259 to R260.
$ws.Rows.Item(194).Insert()

# Populate the new row 194 with the new weekly data point.
$ws.Cells.Item(194, 1).Value = 5
$ws.Cells.Item(194, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(194, 3).Value = "Maule"
$ws.Cells.Item(194, 4).Value = 44559
$ws.Cells.Item(194, 5).Value = 7
$ws.Cells.Item(194, 6).Value = 100112032
$ws.Cells.Item(194, 7).Value = "Zapallo italiano"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Primera"
$ws.Cells.Item(194, 10).Value = 400
$ws.Cells.Item(194, 11).Value = 4000
$ws.Cells.Item(194, 12).Value = 4000
$ws.Cells.Item(194, 13).Value = 4000
$ws.Cells.Item(194, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(194, 15).Value = "Región del Maule"
$ws.Cells.Item(194, 16).Value = 67
$ws.Cells.Item(194, 17).Value = 60
$ws.Cells.Item(194, 18).Value = "Hortaliza"
